{"js": "const oldText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed B\u00fdka 2022: 16.\u201325. ledna\";\nconst newText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16.\u201325. ledna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed B\u00fdka.16.\u201325. ledna\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed B\u00fdka 2022: 16.\u201325. ledna\"\n$newText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 16.\u201325. ledna. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed B\u00fdka.16.\u201325. ledna\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, $null, $null, [ref]$find.Forward, [ref]$find.Wrap, $null, [ref]$find.Replacement.Text, 2)\n"}
